# Bug fixes / tuning of utility-function coefficients in columns M (weight) and N (bias)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = -0.1
$ws.Range("N6").Value = 0.5
$ws.Range("M8").Value = -0.7
$ws.Range("M12").Value = -0.3
$ws.Range("N12").Value = 0.28
$ws.Range("M14").Value = -0.7
$ws.Range("M15").Value = -0.7
$ws.Range("M16").Value = -0.7
$ws.Range("N20").Value = -0.02
$ws.Range("M22").Value = -0.7
$ws.Range("M23").Value = -0.7
$ws.Range("M24").Value = -0.7
$ws.Range("M31").Value = -0.1
$ws.Range("N31").Value = -0.02
$ws.Range("M33").Value = -0.7
$ws.Range("M37").Value = -0.3
$ws.Range("N37").Value = -0.22
$ws.Range("M39").Value = -0.7
$ws.Range("M40").Value = -0.7
$ws.Range("M41").Value = -0.7
$ws.Range("N45").Value = -0.32
$ws.Range("M47").Value = -0.7
$ws.Range("M48").Value = -0.7
$ws.Range("M49").Value = -0.7
$ws.Range("M56").Value = -0.1
$ws.Range("N56").Value = -0.32
$ws.Range("M58").Value = -0.7
$ws.Range("M62").Value = -0.3
$ws.Range("N62").Value = -0.4
$ws.Range("M64").Value = -0.7
$ws.Range("M65").Value = -0.7
$ws.Range("M66").Value = -0.7
$ws.Range("N70").Value = -0.44
$ws.Range("M72").Value = -0.7
$ws.Range("M73").Value = -0.7
$ws.Range("M74").Value = -0.7
$ws.Range("M83").Value = -0.1
$ws.Range("N83").Value = -0.32
$ws.Range("M85").Value = -0.7
$ws.Range("M89").Value = -0.3
$ws.Range("N89").Value = -0.4
$ws.Range("M91").Value = -0.7
$ws.Range("M92").Value = -0.7
$ws.Range("M93").Value = -0.7
$ws.Range("N97").Value = -0.44
$ws.Range("M99").Value = -0.7
$ws.Range("M100").Value = -0.7
$ws.Range("M101").Value = -0.7
$ws.Range("M108").Value = -0.1
$ws.Range("N108").Value = -0.48
$ws.Range("M110").Value = -0.7
$ws.Range("M114").Value = -0.3
$ws.Range("N114").Value = -0.48
$ws.Range("M116").Value = -0.7
$ws.Range("M117").Value = -0.7
$ws.Range("M118").Value = -0.7
$ws.Range("N122").Value = -0.5
$ws.Range("M124").Value = -0.7
$ws.Range("M125").Value = -0.7
$ws.Range("M126").Value = -0.7

# Restore view state: scrolled down to row 94, active cell M104
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M104").Select()
